$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "MaternalDeathWasPregnant" column (AE)
# to host SurgeryWasPerformed and MannerOfDeath.
$ws.Range("AE1:AF1").EntireColumn.Insert()

# Insert one more new column before what is now "MaternalDeathPregnancyContribute" (now at AH)
# to host MaternalDeathTimeFromPregnancy, right after MaternalDeathWasPregnant.
$ws.Range("AH1").EntireColumn.Insert()

# Set the new header labels.
$ws.Range("AE1").Value = "SurgeryWasPerformed"
$ws.Range("AF1").Value = "MannerOfDeath"
$ws.Range("AH1").Value = "MaternalDeathTimeFromPregnancy"

# The old MaternalDeathWasPregnant / MaternalDeathPregnancyContribute sample data
# (0/0 per row) shifted along with the column inserts but is no longer present in
# the updated sample - clear it out.
$ws.Range("AG2:AG5").ClearContents()
$ws.Range("AI2:AI5").ClearContents()

# Update the sample "ICDMinorVersion"-ish D column values for rows 2-5.
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
